# semana 18 de 2025
# Adds a new weekly column (U, header "18") and a new health-provider row
# (6600103414 / 01), shifting the three trailing rows down by one. Also
# corrects T32 (6 -> 3) and fills in T40 (was blank -> 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Structural changes -------------------------------------------------
# Insert the new "18" week column at U (col 21); existing data (cols A:T)
# is untouched.
$ws.Columns.Item(21).Insert()

# Insert a new row at 51; rows 51-53 (EPMSC PEREIRA, SANIDAD POLICIA,
# BATALLON SAN MATEO) shift down to 52-54.
$ws.Rows.Item(51).Insert()

# --- 2. New column U header --------------------------------------------
$ws.Range("U1").Value2 = "'18"

# --- 3. New column U values for existing rows (2-50) -----------------------
$uZeroRows = @(2,4,5,6,7,8,10,11,12,13,17,20,22,23,26,29,30,32,33,34,35,37,38,39,40,41,42,43,44,45,46,47,48,49,50)
foreach ($r in $uZeroRows) {
    $ws.Range("U$r").Value2 = 0
}
$ws.Range("U27").Value2 = 2
$ws.Range("U28").Value2 = 4

# --- 4. Corrections to existing column T values -----------------------
$ws.Range("T32").Value2 = 3
$ws.Range("T40").Value2 = 0

# --- 5. New row 51 content (6600103414 / 01) -------------------------------
# Only cod_pre, cod_sub and the new week column are populated.
$ws.Range("A51").Value2 = "'6600103414"
$ws.Range("B51").Value2 = "'01"
$ws.Range("U51").Value2 = 0

# --- 6. Fix up the two cells that gained data after the shift --------------
# Row 52 (EPMSC PEREIRA) now also reports weeks 9 and 15.
$ws.Range("L52").Value2 = 0
$ws.Range("R52").Value2 = 0
$ws.Range("U52").Value2 = 0

# Row 53 (SANIDAD POLICIA NACIONAL RISARALDA) loses weeks 3, 9 and 15
# (no report that week), so clear those cells back out.
$ws.Range("F53").ClearContents()
$ws.Range("L53").ClearContents()
$ws.Range("R53").ClearContents()
$ws.Range("U53").Value2 = 0

# Row 54 (BATALLON SAN MATEO) keeps its original full set of zeros; just
# stamp the new week column.
$ws.Range("U54").Value2 = 0
